$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Betrayal Legacy" entry (row 61) from the BGG id lookup table.
$ws.Rows.Item(61).Delete()

# Insert the new "Chocolate Factory: Deluxe Edition" entry directly after
# "Chocolate Factory" (now row 98, having shifted up after the delete above),
# keeping the table in alphabetical order.
$ws.Rows.Item(99).Insert()
$ws.Range("A99").Value = "Chocolate Factory: Deluxe Edition"
$ws.Range("B99").Value = 329434
